$wb = $excel.ActiveWorkbook

# --- "dataset" sheet: fix the SUM formula in G2 to include column E as well ---
$dataset = $wb.Worksheets.Item("dataset")
$dataset.Range("G2").Formula = "=SUM(E2:F2)"

# --- "rq2" sheet: update the supporting (helper) figures in columns F:I ---
$rq2 = $wb.Worksheets.Item("rq2")
$rq2.Range("F3").Value = 661
$rq2.Range("G3").Value = 37

$rq2.Range("F4").Value = 653
$rq2.Range("G4").Value = 35
$rq2.Range("H4").Value = 76
$rq2.Range("I4").Value = 334

$rq2.Range("F5").Value = 650
$rq2.Range("G5").Value = 33
$rq2.Range("H5").Value = 82
$rq2.Range("I5").Value = 333

# Refs/print setup for rq2
$rq2.PageSetup.PaperSize = 9
$rq2.PageSetup.Orientation = 1

# --- Switch the active tab from "rq2" to "dataset", updating each sheet's selection ---
$rq2.Range("F6").Select()

$dataset.Activate()
$dataset.Range("J6").Select()
